$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1282474195497992, 0.9763931104222339, 4.446347975453152, 2.108636520468417, 2.125676030426229, 51)
    3  = @(0.1215397788006174, 1.018530036818849, 4.443793105881769, 2.108030622614807, 2.125890254825519, 50)
    4  = @(0.1466579445676272, 0.9414795404354654, 4.22668197375337, 2.055889582091745, 2.071902797105537, 49)
    5  = @(0.1692583036833795, 1.044790455529975, 4.665004838604097, 2.159862226764498, 2.176006054099445, 48)
    6  = @(0.1384867545238795, 0.9904347050370882, 4.525272260743743, 2.127268732610843, 2.145705533808543, 47)
    7  = @(0.1625861655212504, 1.003062538999909, 4.667669780972878, 2.160479062840665, 2.17815841638289, 46)
    8  = @(0.09841386822182357, 0.9352685283506139, 4.357710809747646, 2.087513068162124, 2.108754201445493, 45)
    9  = @(0.06748005575673845, 0.9265491817941878, 4.470640045755423, 2.114388811395724, 2.137743895861163, 44)
    10 = @(0.09914051671729636, 0.9226436494409858, 4.481900347334963, 2.117049916117937, 2.13975451822466, 43)
    11 = @(0.07576832823998672, 0.9176283451160873, 4.557762439856401, 2.13489166934915, 2.159408835161238, 42)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
